$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension-affecting data: rows 2-13 get reindexed string values (Sema5a/Plxnb3 shared string ids shift)
# and numeric recalculated values; rows 14-17 are brand new rows for the "Resolving-Mac" sending cluster.

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Sema5a"
$ws.Range("C2").Value = "Plxnb3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.4993116666666666
$ws.Range("H2").Value = 1.497935
$ws.Range("I2").Value = 0.01118655668236004
$ws.Range("J2").Value = 0.01118655668236004
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4640010000000001
$ws.Range("N2").Value = 1.392003
$ws.Range("O2").Value = 0.128664373327748
$ws.Range("P2").Value = 0.128664373327748
$ws.Range("Q2").Value = 0.231681112645
$ws.Range("R2").Value = 2.085130013805
$ws.Range("S2").Value = 0.001439311305231187
$ws.Range("T2").Value = 0.001439311305231187

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Sema5a"
$ws.Range("C3").Value = "Plxnb3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.4993116666666666
$ws.Range("H3").Value = 1.497935
$ws.Range("I3").Value = 0.01118655668236004
$ws.Range("J3").Value = 0.01118655668236004
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.5337323333333334
$ws.Range("N3").Value = 1.601197
$ws.Range("O3").Value = 0.148000405587682
$ws.Range("P3").Value = 0.148000405587682
$ws.Range("Q3").Value = 0.2664987809105555
$ws.Range("R3").Value = 2.398489028195
$ws.Range("S3").Value = 0.00165561492611888
$ws.Range("T3").Value = 0.00165561492611888

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Sema5a"
$ws.Range("C4").Value = "Plxnb3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.4993116666666666
$ws.Range("H4").Value = 1.497935
$ws.Range("I4").Value = 0.01118655668236004
$ws.Range("J4").Value = 0.01118655668236004
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.709312
$ws.Range("N4").Value = 5.127936
$ws.Range("O4").Value = 0.4739807830190013
$ws.Range("P4").Value = 0.4739807830190014
$ws.Range("Q4").Value = 0.8534794235733332
$ws.Range("R4").Value = 7.681314812159999
$ws.Range("S4").Value = 0.005302212895591452
$ws.Range("T4").Value = 0.005302212895591453

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Sema5a"
$ws.Range("C5").Value = "Plxnb3"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.4993116666666666
$ws.Range("H5").Value = 1.497935
$ws.Range("I5").Value = 0.01118655668236004
$ws.Range("J5").Value = 0.01118655668236004
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.8992443333333332
$ws.Range("N5").Value = 2.697732999999999
$ws.Range("O5").Value = 0.2493544380655686
$ws.Range("P5").Value = 0.2493544380655686
$ws.Range("Q5").Value = 0.4490031868172221
$ws.Range("R5").Value = 4.041028681354999
$ws.Range("S5").Value = 0.002789417555418518
$ws.Range("T5").Value = 0.002789417555418518

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Sema5a"
$ws.Range("C6").Value = "Plxnb3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 32.04971933333334
$ws.Range("H6").Value = 96.149158
$ws.Range("I6").Value = 0.718040506382581
$ws.Range("J6").Value = 0.7180405063825809
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.4640010000000001
$ws.Range("N6").Value = 1.392003
$ws.Range("O6").Value = 0.128664373327748
$ws.Range("P6").Value = 0.128664373327748
$ws.Range("Q6").Value = 14.871101820386
$ws.Range("R6").Value = 133.839916383474
$ws.Range("S6").Value = 0.09238623177765366
$ws.Range("T6").Value = 0.09238623177765365

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Sema5a"
$ws.Range("C7").Value = "Plxnb3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 32.04971933333334
$ws.Range("H7").Value = 96.149158
$ws.Range("I7").Value = 0.718040506382581
$ws.Range("J7").Value = 0.7180405063825809
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.5337323333333334
$ws.Range("N7").Value = 1.601197
$ws.Range("O7").Value = 0.148000405587682
$ws.Range("P7").Value = 0.148000405587682
$ws.Range("Q7").Value = 17.10597148245845
$ws.Range("R7").Value = 153.953743342126
$ws.Range("S7").Value = 0.1062702861730066
$ws.Range("T7").Value = 0.1062702861730066

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Sema5a"
$ws.Range("C8").Value = "Plxnb3"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 32.04971933333334
$ws.Range("H8").Value = 96.149158
$ws.Range("I8").Value = 0.718040506382581
$ws.Range("J8").Value = 0.7180405063825809
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.709312
$ws.Range("N8").Value = 5.127936
$ws.Range("O8").Value = 0.4739807830190013
$ws.Range("P8").Value = 0.4739807830190014
$ws.Range("Q8").Value = 54.78296985309867
$ws.Range("R8").Value = 493.046728677888
$ws.Range("S8").Value = 0.340337401454576
$ws.Range("T8").Value = 0.3403374014545759

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Sema5a"
$ws.Range("C9").Value = "Plxnb3"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 32.04971933333334
$ws.Range("H9").Value = 96.149158
$ws.Range("I9").Value = 0.718040506382581
$ws.Range("J9").Value = 0.7180405063825809
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.8992443333333332
$ws.Range("N9").Value = 2.697732999999999
$ws.Range("O9").Value = 0.2493544380655686
$ws.Range("P9").Value = 0.2493544380655686
$ws.Range("Q9").Value = 28.82052849542378
$ws.Range("R9").Value = 259.384756458814
$ws.Range("S9").Value = 0.1790465869773448
$ws.Range("T9").Value = 0.1790465869773447

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Sema5a"
$ws.Range("C10").Value = "Plxnb3"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 12.08377133333333
$ws.Range("H10").Value = 36.251314
$ws.Range("I10").Value = 0.2707242830102989
$ws.Range("J10").Value = 0.2707242830102989
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.4640010000000001
$ws.Range("N10").Value = 1.392003
$ws.Range("O10").Value = 0.128664373327748
$ws.Range("P10").Value = 0.128664373327748
$ws.Range("Q10").Value = 5.606881982438001
$ws.Range("R10").Value = 50.461937841942
$ws.Range("S10").Value = 0.03483257021812401
$ws.Range("T10").Value = 0.03483257021812402

# Row 11
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Sema5a"
$ws.Range("C11").Value = "Plxnb3"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 12.08377133333333
$ws.Range("H11").Value = 36.251314
$ws.Range("I11").Value = 0.2707242830102989
$ws.Range("J11").Value = 0.2707242830102989
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.5337323333333334
$ws.Range("N11").Value = 1.601197
$ws.Range("O11").Value = 0.148000405587682
$ws.Range("P11").Value = 0.148000405587682
$ws.Range("Q11").Value = 6.449499469206445
$ws.Range("R11").Value = 58.045495222858
$ws.Range("S11").Value = 0.04006730368795865
$ws.Range("T11").Value = 0.04006730368795866

# Row 12
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Sema5a"
$ws.Range("C12").Value = "Plxnb3"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 12.08377133333333
$ws.Range("H12").Value = 36.251314
$ws.Range("I12").Value = 0.2707242830102989
$ws.Range("J12").Value = 0.2707242830102989
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 1.709312
$ws.Range("N12").Value = 5.127936
$ws.Range("O12").Value = 0.4739807830190013
$ws.Range("P12").Value = 0.4739807830190014
$ws.Range("Q12").Value = 20.65493534532267
$ws.Range("R12").Value = 185.894418107904
$ws.Range("S12").Value = 0.1283181076434792
$ws.Range("T12").Value = 0.1283181076434792

# Row 13
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Sema5a"
$ws.Range("C13").Value = "Plxnb3"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 12.08377133333333
$ws.Range("H13").Value = 36.251314
$ws.Range("I13").Value = 0.2707242830102989
$ws.Range("J13").Value = 0.2707242830102989
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.8992443333333332
$ws.Range("N13").Value = 2.697732999999999
$ws.Range("O13").Value = 0.2493544380655686
$ws.Range("P13").Value = 0.2493544380655686
$ws.Range("Q13").Value = 10.86626289679578
$ws.Range("R13").Value = 97.79636607116198
$ws.Range("S13").Value = 0.06750630146073702
$ws.Range("T13").Value = 0.06750630146073704

# Row 14
$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "Sema5a"
$ws.Range("C14").Value = "Plxnb3"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.002171666666666667
$ws.Range("H14").Value = 0.006515
$ws.Range("I14").Value = ([double]"4.86539247601369e-05")
$ws.Range("J14").Value = ([double]"4.86539247601369e-05")
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.4640010000000001
$ws.Range("N14").Value = 1.392003
$ws.Range("O14").Value = 0.128664373327748
$ws.Range("P14").Value = 0.128664373327748
$ws.Range("Q14").Value = 0.001007655505
$ws.Range("R14").Value = 0.009068899545
$ws.Range("S14").Value = ([double]"6.260026739198418e-06")
$ws.Range("T14").Value = ([double]"6.260026739198418e-06")

# Row 15
$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "Sema5a"
$ws.Range("C15").Value = "Plxnb3"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.002171666666666667
$ws.Range("H15").Value = 0.006515
$ws.Range("I15").Value = ([double]"4.86539247601369e-05")
$ws.Range("J15").Value = ([double]"4.86539247601369e-05")
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0.5337323333333334
$ws.Range("N15").Value = 1.601197
$ws.Range("O15").Value = 0.148000405587682
$ws.Range("P15").Value = 0.148000405587682
$ws.Range("Q15").Value = 0.001159088717222222
$ws.Range("R15").Value = 0.010431798455
$ws.Range("S15").Value = ([double]"7.200800597932826e-06")
$ws.Range("T15").Value = ([double]"7.200800597932826e-06")

# Row 16
$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "Sema5a"
$ws.Range("C16").Value = "Plxnb3"
$ws.Range("D16").Value = "MuSCs"
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.002171666666666667
$ws.Range("H16").Value = 0.006515
$ws.Range("I16").Value = ([double]"4.86539247601369e-05")
$ws.Range("J16").Value = ([double]"4.86539247601369e-05")
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1.709312
$ws.Range("N16").Value = 5.127936
$ws.Range("O16").Value = 0.4739807830190013
$ws.Range("P16").Value = 0.4739807830190014
$ws.Range("Q16").Value = 0.003712055893333333
$ws.Range("R16").Value = 0.03340850304
$ws.Range("S16").Value = ([double]"2.306102535475726e-05")
$ws.Range("T16").Value = ([double]"2.306102535475727e-05")

# Row 17
$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("B17").Value = "Sema5a"
$ws.Range("C17").Value = "Plxnb3"
$ws.Range("D17").Value = "Resolving-Mac"
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.002171666666666667
$ws.Range("H17").Value = 0.006515
$ws.Range("I17").Value = ([double]"4.86539247601369e-05")
$ws.Range("J17").Value = ([double]"4.86539247601369e-05")
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.8992443333333332
$ws.Range("N17").Value = 2.697732999999999
$ws.Range("O17").Value = 0.2493544380655686
$ws.Range("P17").Value = 0.2493544380655686
$ws.Range("Q17").Value = 0.001952858943888889
$ws.Range("R17").Value = 0.017575730495
$ws.Range("S17").Value = ([double]"1.213207206824839e-05")
$ws.Range("T17").Value = ([double]"1.213207206824839e-05")
